$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.802.23"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.628.69"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.28"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5076"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2579"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06437"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07799"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.261"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.624.75"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5582"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.15"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₅7566"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.811.11"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.68"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.310"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.811"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.007"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.797"
$ws.Range("E25").Value = "  -4.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.16"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1257"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.734"
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.39"
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.236"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04872"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.277"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.190"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.556"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.373"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8952"
$ws.Range("E36").Value = "  -2.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.568"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.132.78"
$ws.Range("E38").Value = "  +3.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5468"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01557"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9982"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.561"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7957"
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.30"
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.780.78"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -7.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4430"
$ws.Range("E47").Value = "  -2.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.10"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05051"
$ws.Range("E49").Value = "  -2.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.604"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  +0.20%  "
